$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add three new instrument rows (Bass Clarinet, Cornet, C Tenor Saxophone).
# Column order matches how the shared-string table ends up laid out:
# English names first, then the internal "key" codes, then the German names.

$ws.Range("B11").Value = "Bass Clarinet"
$ws.Range("B12").Value = "Cornet"
$ws.Range("B13").Value = "C Tenor Saxophone"

$ws.Range("A11").Value = "bass_clarinet"
$ws.Range("A12").Value = "cornet"
$ws.Range("A13").Value = "c_tenor_saxophone"

$ws.Range("C11").Value = "Bassklarinette"
$ws.Range("C12").Value = "Kornett"
$ws.Range("C13").Value = "C-Tenor Saxophon"

$ws.Range("D11").Value = 34
$ws.Range("E11").Value = 83
$ws.Range("F11").Value = 2
$ws.Range("G11").Value = "treble"

$ws.Range("D12").Value = 54
$ws.Range("E12").Value = 84
$ws.Range("F12").Value = 2
$ws.Range("G12").Value = "treble"

$ws.Range("D13").Value = 44
$ws.Range("E13").Value = 76
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = "treble"

# Move the active selection like the author left it.
$ws.Range("I12").Select()
